$wb = $excel.ActiveWorkbook

$oldVersion = "mines - version 1.0.0 (Feb 3 2026) (built on February 03 2026 10.14.00 EST)"
$newVersion = "Coal Mine Boundaries and Methane Sources - version 1.0.0 (built on February 03 2026 17.29.55 EST)"

$aboutSheet = $wb.Worksheets.Item("About")
$dataSheet = $wb.Worksheets.Item("Boundaries and methane sources")

# Update the "About" sheet version line (A2)
$aboutSheet.Range("A2").Value = "Version: $newVersion"

# Update the "About" sheet recommended citation line (A6)
$newCitation = "Recommended Citation:  `"Global Energy Monitor, Coal mine boundaries and methane sources for Kirova Coal Mine, Russia, M0809, version '$newVersion'. (See the CC license for attribution requirements if sharing or adapting the data set.)"
$aboutSheet.Range("A6").Value = $newCitation

# Update the build_version column (S) on the data sheet for all data rows (2-36)
$usedRange = $dataSheet.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $dataSheet.Cells.Item($r, 19)  # Column S = 19
    if ($cell.Value2 -eq $oldVersion) {
        $cell.Value2 = $newVersion
    }
}
